# Apply "Minor Changes to general.xlsx (translations)":
# Add a new translation key "of the World" / "der Welt" as the next row
# (row 55) on the Exiobase and Deutsch sheets, and row 59 on the English
# sheet, then leave the Exiobase sheet as the active tab/selection.

$wb = $excel.ActiveWorkbook

# --- Exiobase sheet: new row 55, English key used in both columns ---
$wsExio = $wb.Worksheets.Item("Exiobase")
$wsExio.Range("A55").Value = "of the World"
$wsExio.Range("B55").Value = "of the World"
$wsExio.Rows.Item(55).RowHeight = 15

# --- Deutsch sheet: new row 55, English key + German translation ---
$wsDe = $wb.Worksheets.Item("Deutsch")
$wsDe.Range("A55").Value = "of the World"
$wsDe.Range("B55").Value = "der Welt"
$wsDe.Rows.Item(55).RowHeight = 15

# --- English sheet: new row 59, English key used in both columns ---
$wsEn = $wb.Worksheets.Item("English")
$wsEn.Range("A59").Value = "of the World"
$wsEn.Range("B59").Value = "of the World"
$wsEn.Rows.Item(59).RowHeight = 15

# --- Final view/selection state ---
# English sheet ends up with a plain (non-tab-selected) selection at B59
[void]$wsEn.Select()
[void]$wsEn.Range("B59").Select()

# Deutsch sheet ends up selected at A55
[void]$wsDe.Select()
[void]$wsDe.Range("A55").Select()

# Exiobase ends up the active/selected tab with selection at B57
[void]$wsExio.Select()
[void]$wsExio.Range("B57").Select()
